$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 24 de Abril de 2020 a las 18:22'
$ws.Cells.Item(4, 2).Value = 894034
$ws.Cells.Item(4, 3).Value = 7592
$ws.Cells.Item(4, 4).Value = 90187
$ws.Cells.Item(4, 5).Value = 752928
$ws.Cells.Item(4, 7).Value = 683
$ws.Cells.Item(4, 8).Value = 50919
$ws.Cells.Item(16, 2).Value = 42773
$ws.Cells.Item(16, 3).Value = 663
$ws.Cells.Item(16, 4).Value = 15202
$ws.Cells.Item(16, 5).Value = 25374
$ws.Cells.Item(16, 7).Value = 50
$ws.Cells.Item(16, 8).Value = 2197
$ws.Cells.Item(28, 1).Value = 'Chile'
$ws.Cells.Item(28, 2).Value = 12306
$ws.Cells.Item(28, 3).Value = 494
$ws.Cells.Item(28, 4).Value = 6327
$ws.Cells.Item(28, 5).Value = 5805
$ws.Cells.Item(28, 6).Value = 408
$ws.Cells.Item(28, 7).Value = 6
$ws.Cells.Item(28, 8).Value = 174
$ws.Cells.Item(29, 1).Value = 'Singapur'
$ws.Cells.Item(29, 2).Value = 12075
$ws.Cells.Item(29, 3).Value = 897
$ws.Cells.Item(29, 4).Value = 924
$ws.Cells.Item(29, 5).Value = 11139
$ws.Cells.Item(29, 6).Value = 26
$ws.Cells.Item(29, 8).Value = 12
$ws.Cells.Item(33, 2).Value = 10892
$ws.Cells.Item(33, 3).Value = 381
$ws.Cells.Item(33, 5).Value = 8454
$ws.Cells.Item(33, 7).Value = 40
$ws.Cells.Item(33, 8).Value = 494
$ws.Cells.Item(47, 5).Value = 4719
$ws.Cells.Item(47, 7).Value = 2
$ws.Cells.Item(47, 8).Value = 267
$ws.Cells.Item(82, 1).Value = 'Cuba'
$ws.Cells.Item(82, 2).Value = 1285
$ws.Cells.Item(82, 3).Value = 50
$ws.Cells.Item(82, 4).Value = 416
$ws.Cells.Item(82, 5).Value = 820
$ws.Cells.Item(82, 6).Value = 11
$ws.Cells.Item(82, 7).Value = 6
$ws.Cells.Item(82, 8).Value = 49
$ws.Cells.Item(83, 1).Value = 'Ghana'
$ws.Cells.Item(83, 2).Value = 1279
$ws.Cells.Item(83, 3).Value = 125
$ws.Cells.Item(83, 4).Value = 134
$ws.Cells.Item(83, 5).Value = 1135
$ws.Cells.Item(83, 6).Value = 4
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 10
$ws.Cells.Item(93, 2).Value = 731
$ws.Cells.Item(93, 3).Value = 8
$ws.Cells.Item(93, 5).Value = 358
$ws.Cells.Item(93, 7).Value = 3
$ws.Cells.Item(93, 8).Value = 40
$ws.Cells.Item(105, 2).Value = 513
$ws.Cells.Item(105, 3).Value = 12
$ws.Cells.Item(105, 4).Value = 64
$ws.Cells.Item(105, 5).Value = 409
$ws.Cells.Item(108, 2).Value = 441
$ws.Cells.Item(108, 3).Value = 4
$ws.Cells.Item(108, 4).Value = 326
$ws.Cells.Item(108, 5).Value = 108
$ws.Cells.Item(117, 4).Value = 285
$ws.Cells.Item(117, 5).Value = 37
$ws.Cells.Item(164, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(164, 3).Value = 49
$ws.Cells.Item(164, 4).Value = 3
$ws.Cells.Item(164, 5).Value = 48
$ws.Cells.Item(164, 7).Value = 2
$ws.Cells.Item(164, 8).Value = 3
$ws.Cells.Item(165, 1).Value = 'Benin'
$ws.Cells.Item(165, 2).Value = 54
$ws.Cells.Item(165, 4).Value = 27
$ws.Cells.Item(165, 5).Value = 26
$ws.Cells.Item(165, 8).Value = 1
$ws.Cells.Item(166, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(166, 2).Value = 50
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 3
$ws.Cells.Item(166, 5).Value = 47
$ws.Cells.Item(167, 1).Value = 'Nepal'
$ws.Cells.Item(167, 2).Value = 49
$ws.Cells.Item(167, 3).Value = 1
$ws.Cells.Item(167, 4).Value = 10
$ws.Cells.Item(167, 5).Value = 39
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(168, 1).Value = 'Macao'
$ws.Cells.Item(168, 2).Value = 45
$ws.Cells.Item(168, 4).Value = 27
$ws.Cells.Item(168, 5).Value = 18
$ws.Cells.Item(168, 6).Value = 1
$ws.Cells.Item(168, 8).Value = 0
$ws.Cells.Item(169, 1).Value = 'Siria'
$ws.Cells.Item(169, 2).Value = 42
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(169, 4).Value = 6
$ws.Cells.Item(169, 5).Value = 33
$ws.Cells.Item(169, 8).Value = 3
$ws.Cells.Item(170, 1).Value = 'Republica del Chad'
$ws.Cells.Item(170, 2).Value = 40
$ws.Cells.Item(170, 3).Value = 7
$ws.Cells.Item(170, 4).Value = 8
$ws.Cells.Item(170, 5).Value = 32
$ws.Cells.Item(170, 8).Value = 0
$ws.Cells.Item(171, 1).Value = 'Puerto Rico'
$ws.Cells.Item(171, 4).Value = 1
$ws.Cells.Item(171, 5).Value = 36
$ws.Cells.Item(171, 8).Value = 2
$ws.Cells.Item(172, 1).Value = 'Eritrea'
$ws.Cells.Item(172, 2).Value = 39
$ws.Cells.Item(172, 4).Value = 11
$ws.Cells.Item(172, 5).Value = 28
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 8).Value = 0
$ws.Cells.Item(173, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(173, 2).Value = 38
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 20
$ws.Cells.Item(173, 5).Value = 16
$ws.Cells.Item(173, 6).Value = 3
$ws.Cells.Item(173, 8).Value = 2
$ws.Cells.Item(174, 1).Value = 'Mongolia'
$ws.Cells.Item(174, 3).Value = 1
$ws.Cells.Item(174, 4).Value = 9
$ws.Cells.Item(174, 5).Value = 27
$ws.Cells.Item(174, 8).Value = 0
$ws.Cells.Item(175, 1).Value = 'Suazilandia'
$ws.Cells.Item(175, 2).Value = 36
$ws.Cells.Item(175, 3).Value = 5
$ws.Cells.Item(175, 4).Value = 10
$ws.Cells.Item(175, 5).Value = 25
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 8).Value = 1
$ws.Cells.Item(176, 1).Value = 'Malaui'
$ws.Cells.Item(176, 2).Value = 33
$ws.Cells.Item(176, 4).Value = 3
$ws.Cells.Item(176, 5).Value = 27
$ws.Cells.Item(176, 6).Value = 1
$ws.Cells.Item(176, 8).Value = 3
$ws.Cells.Item(177, 1).Value = 'Guam'
$ws.Cells.Item(177, 2).Value = 32
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 31
$ws.Cells.Item(177, 8).Value = 1
$ws.Cells.Item(178, 1).Value = 'Zimbabue'
$ws.Cells.Item(178, 2).Value = 29
$ws.Cells.Item(178, 3).Value = 1
$ws.Cells.Item(178, 4).Value = 2
$ws.Cells.Item(178, 5).Value = 23
$ws.Cells.Item(178, 8).Value = 4
$ws.Cells.Item(179, 1).Value = 'Angola'
$ws.Cells.Item(179, 2).Value = 25
$ws.Cells.Item(179, 4).Value = 6
$ws.Cells.Item(179, 5).Value = 17
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 8).Value = 2
$ws.Cells.Item(180, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(180, 2).Value = 24
$ws.Cells.Item(180, 4).Value = 10
$ws.Cells.Item(180, 5).Value = 11
$ws.Cells.Item(180, 6).Value = 1
$ws.Cells.Item(180, 8).Value = 3
$ws.Cells.Item(181, 1).Value = 'Timor Oriental'
$ws.Cells.Item(181, 2).Value = 23
$ws.Cells.Item(181, 4).Value = 1
$ws.Cells.Item(181, 5).Value = 22
$ws.Cells.Item(181, 8).Value = 0
$ws.Cells.Item(182, 1).Value = 'Botsuana'
$ws.Cells.Item(182, 2).Value = 22
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 21
$ws.Cells.Item(182, 8).Value = 1
$ws.Cells.Item(183, 1).Value = 'Laos'
$ws.Cells.Item(183, 2).Value = 19
$ws.Cells.Item(183, 4).Value = 4
$ws.Cells.Item(183, 5).Value = 15
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 8).Value = 0
$ws.Cells.Item(184, 1).Value = 'Belice'
$ws.Cells.Item(184, 4).Value = 5
$ws.Cells.Item(184, 5).Value = 11
$ws.Cells.Item(184, 6).Value = 1
$ws.Cells.Item(184, 8).Value = 2
$ws.Cells.Item(185, 1).Value = 'Fiyi'
$ws.Cells.Item(185, 4).Value = 10
$ws.Cells.Item(185, 5).Value = 8
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(186, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(186, 2).Value = 18
$ws.Cells.Item(186, 4).Value = 17
$ws.Cells.Item(186, 5).Value = 1
$ws.Cells.Item(186, 6).Value = 1
$ws.Cells.Item(187, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(187, 2).Value = 17
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 17
$ws.Cells.Item(188, 1).Value = 'Namibia'
$ws.Cells.Item(188, 4).Value = 7
$ws.Cells.Item(188, 5).Value = 9
$ws.Cells.Item(189, 1).Value = 'Dominica'
$ws.Cells.Item(189, 4).Value = 9
$ws.Cells.Item(189, 5).Value = 7
$ws.Cells.Item(190, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(190, 2).Value = 16
$ws.Cells.Item(190, 4).Value = 10
$ws.Cells.Item(190, 5).Value = 6
$ws.Cells.Item(191, 1).Value = 'San Cristobal y Nieves'
$ws.Cells.Item(191, 4).Value = 2
$ws.Cells.Item(191, 5).Value = 13
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(192, 1).Value = 'Granada'
$ws.Cells.Item(192, 4).Value = 7
$ws.Cells.Item(192, 5).Value = 8
$ws.Cells.Item(192, 6).Value = 4
$ws.Cells.Item(193, 1).Value = 'Santa Lucia'
$ws.Cells.Item(193, 2).Value = 15
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 15
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(194, 1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(194, 3).Value = 1
$ws.Cells.Item(194, 4).Value = 5
$ws.Cells.Item(194, 5).Value = 9
$ws.Cells.Item(194, 8).Value = 0
$ws.Cells.Item(195, 1).Value = 'Curazao'
$ws.Cells.Item(195, 2).Value = 14
$ws.Cells.Item(195, 5).Value = 2
$ws.Cells.Item(195, 8).Value = 1
$ws.Cells.Item(196, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(196, 2).Value = 12
$ws.Cells.Item(196, 4).Value = 11
$ws.Cells.Item(196, 5).Value = 1
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(197, 1).Value = 'Montserrat'
$ws.Cells.Item(197, 4).Value = 2
$ws.Cells.Item(197, 5).Value = 9
$ws.Cells.Item(197, 6).Value = 1
$ws.Cells.Item(197, 8).Value = 0
$ws.Cells.Item(198, 1).Value = 'Burundi'
$ws.Cells.Item(199, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(199, 4).Value = 4
$ws.Cells.Item(199, 5).Value = 6
$ws.Cells.Item(199, 8).Value = 1
$ws.Cells.Item(200, 1).Value = 'Seychelles'
$ws.Cells.Item(200, 4).Value = 6
$ws.Cells.Item(200, 5).Value = 5
$ws.Cells.Item(200, 8).Value = 0
$ws.Cells.Item(201, 1).Value = 'Nicaragua'
$ws.Cells.Item(201, 4).Value = 7
$ws.Cells.Item(201, 5).Value = 1
$ws.Cells.Item(201, 8).Value = 3
$ws.Cells.Item(202, 1).Value = 'Groenlandia'
$ws.Cells.Item(202, 2).Value = 11
$ws.Cells.Item(202, 4).Value = 11
$ws.Cells.Item(202, 5).Value = 0
$ws.Cells.Item(202, 8).Value = 0
$ws.Cells.Item(203, 1).Value = 'Gambia'
$ws.Cells.Item(203, 4).Value = 2
$ws.Cells.Item(203, 5).Value = 7
$ws.Cells.Item(204, 1).Value = 'Surinam'
$ws.Cells.Item(204, 2).Value = 10
$ws.Cells.Item(204, 4).Value = 6
$ws.Cells.Item(204, 5).Value = 3
$ws.Cells.Item(204, 8).Value = 1
$ws.Cells.Item(205, 1).Value = 'Santa Sede'
$ws.Cells.Item(205, 2).Value = 9
$ws.Cells.Item(205, 4).Value = 2
$ws.Cells.Item(205, 5).Value = 7
$ws.Cells.Item(206, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(206, 2).Value = 8
$ws.Cells.Item(206, 4).Value = 0
$ws.Cells.Item(206, 5).Value = 8
$ws.Cells.Item(207, 1).Value = 'Butan'
$ws.Cells.Item(207, 4).Value = 3
$ws.Cells.Item(207, 5).Value = 4
$ws.Cells.Item(207, 8).Value = 0
$ws.Cells.Item(208, 1).Value = 'Mauritania'
$ws.Cells.Item(208, 2).Value = 7
$ws.Cells.Item(208, 4).Value = 6
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 8).Value = 1
$ws.Cells.Item(209, 1).Value = 'Sahara Occidental'
$ws.Cells.Item(209, 4).Value = 5
$ws.Cells.Item(209, 5).Value = 1
$ws.Cells.Item(210, 1).Value = 'San Bartolome'
$ws.Cells.Item(210, 2).Value = 6
$ws.Cells.Item(210, 4).Value = 6
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(211, 1).Value = 'Bonaire, San Eustaquio y Saba'
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(212, 1).Value = 'Sudan del Sur'
$ws.Cells.Item(212, 3).Value = 1
$ws.Cells.Item(212, 4).Value = 0
$ws.Cells.Item(212, 5).Value = 5
$ws.Cells.Item(212, 8).Value = 0
